$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.025425
$ws.Range("H2").Value = 0.076275
$ws.Range("I2").Value = 0.5785685028141451
$ws.Range("J2").Value = 0.578568502814145
$ws.Range("Q2").Value = 0.31415314755
$ws.Range("R2").Value = 2.82737832795
$ws.Range("S2").Value = 0.5785685028141451
$ws.Range("T2").Value = 0.578568502814145

# Row 3 updates
$ws.Range("I3").Value = 0.4214314971858549
$ws.Range("J3").Value = 0.4214314971858549
$ws.Range("S3").Value = 0.4214314971858549
$ws.Range("T3").Value = 0.4214314971858549
